$d = $word.ActiveDocument

$replacements = @(
    @{old="241×3=723"; new="382×8=3056"},
    @{old="485×2=970"; new="287×9=2583"},
    @{old="946×3=2838"; new="826×4=3304"},
    @{old="626×3=1878"; new="854×8=6832"},
    @{old="425×6=2550"; new="145×2=290"},
    @{old="307×9=2763"; new="368×5=1840"},
    @{old="616×9=5544"; new="987×5=4935"},
    @{old="215×6=1290"; new="881×3=2643"},
    @{old="555×7=3885"; new="421×2=842"},
    @{old="179×9=1611"; new="729×9=6561"},
    @{old="327×2=654"; new="481×6=2886"},
    @{old="543×4=2172"; new="699×9=6291"},
    @{old="440×6=2640"; new="913×5=4565"},
    @{old="697×8=5576"; new="759×7=5313"},
    @{old="649×9=5841"; new="671×2=1342"},
    @{old="137×4=548"; new="502×8=4016"},
    @{old="893×5=4465"; new="636×9=5724"},
    @{old="669×4=2676"; new="742×2=1484"},
    @{old="231×9=2079"; new="120×6=720"},
    @{old="549×4=2196"; new="993×9=8937"},
    @{old="145×6=870"; new="300×9=2700"},
    @{old="230×9=2070"; new="331×5=1655"},
    @{old="310×6=1860"; new="128×6=768"},
    @{old="323×8=2584"; new="520×6=3120"},
    @{old="625×5=3125"; new="105×4=420"}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
